# Updated cryptos list with GitHub Actions
# Applies per-cell text updates while preserving the original
# "inline string / text" cell semantics (Excel would otherwise silently
# coerce number-looking strings like "0.9972" into numeric cells).

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range("D2") "29.647.46"
Set-TextValue $ws.Range("E2") "  +1.87%  "
Set-TextValue $ws.Range("D3") "1.847.40"
Set-TextValue $ws.Range("E3") "  +1.11%  "
Set-TextValue $ws.Range("D4") "0.9972"
Set-TextValue $ws.Range("E4") "  -1.08%  "
Set-TextValue $ws.Range("D5") "244.53"
Set-TextValue $ws.Range("E5") "  +0.89%  "
Set-TextValue $ws.Range("D6") "0.6319"
Set-TextValue $ws.Range("E6") "  +3.28%  "
Set-TextValue $ws.Range("D7") "0.9989"
Set-TextValue $ws.Range("E7") "  -1.01%  "
Set-TextValue $ws.Range("B8") "Cardano"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue $ws.Range("D8") "0.2987"
Set-TextValue $ws.Range("E8") "  +3.57%  "
Set-TextValue $ws.Range("B9") "Dogecoin"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D9") "0.07472"
Set-TextValue $ws.Range("E9") "  +1.99%  "
Set-TextValue $ws.Range("D10") "23.86"
Set-TextValue $ws.Range("E10") "  +5.07%  "
Set-TextValue $ws.Range("D11") "0.07674"
Set-TextValue $ws.Range("E11") "  -0.53%  "
Set-TextValue $ws.Range("D12") "1.834.22"
Set-TextValue $ws.Range("E12") "  +0.93%  "
Set-TextValue $ws.Range("D13") "5.041"
Set-TextValue $ws.Range("E13") "  +2.25%  "
Set-TextValue $ws.Range("D14") "0.6850"
Set-TextValue $ws.Range("E14") "  +4.15%  "
Set-TextValue $ws.Range("D15") "84.37"
Set-TextValue $ws.Range("E15") "  +3.65%  "
Set-TextValue $ws.Range("D16") "0.000009383"
Set-TextValue $ws.Range("E16") "  +5.60%  "
Set-TextValue $ws.Range("D17") "5.999"
Set-TextValue $ws.Range("E17") "  +2.63%  "
Set-TextValue $ws.Range("D18") "29.591.98"
Set-TextValue $ws.Range("E18") "  +1.79%  "
Set-TextValue $ws.Range("D19") "2.089.30"
Set-TextValue $ws.Range("E19") "  +1.19%  "
Set-TextValue $ws.Range("D20") "239.02"
Set-TextValue $ws.Range("E20") "  +1.47%  "
Set-TextValue $ws.Range("D21") "12.61"
Set-TextValue $ws.Range("E21") "  +1.49%  "
Set-TextValue $ws.Range("D22") "0.9992"
Set-TextValue $ws.Range("E22") "  -1.02%  "
Set-TextValue $ws.Range("D23") "7.363"
Set-TextValue $ws.Range("E23") "  +3.84%  "
Set-TextValue $ws.Range("D24") "0.9999"
Set-TextValue $ws.Range("E24") "  -1.06%  "
Set-TextValue $ws.Range("D25") "158.92"
Set-TextValue $ws.Range("E25") "  -0.19%  "
Set-TextValue $ws.Range("D26") "0.1420"
Set-TextValue $ws.Range("E26") "  +1.82%  "
Set-TextValue $ws.Range("D27") "8.533"
Set-TextValue $ws.Range("E27") "  +1.46%  "
Set-TextValue $ws.Range("D28") "17.85"
Set-TextValue $ws.Range("E28") "  +1.37%  "
Set-TextValue $ws.Range("D29") "1.500"
Set-TextValue $ws.Range("E29") "  +0.79%  "
Set-TextValue $ws.Range("D30") "0.06014"
Set-TextValue $ws.Range("E30") "  +7.90%  "
Set-TextValue $ws.Range("D31") "1.252"
Set-TextValue $ws.Range("E31") "  +3.53%  "
Set-TextValue $ws.Range("D32") "4.138"
Set-TextValue $ws.Range("E32") "  +1.69%  "
Set-TextValue $ws.Range("D33") "4.127"
Set-TextValue $ws.Range("E33") "  +1.25%  "
Set-TextValue $ws.Range("D34") "1.869"
Set-TextValue $ws.Range("E34") "  +2.27%  "
Set-TextValue $ws.Range("D35") "1.148"
Set-TextValue $ws.Range("D36") "0.7289"
Set-TextValue $ws.Range("E36") "  -0.35%  "
Set-TextValue $ws.Range("D37") "2.610"
Set-TextValue $ws.Range("E37") "  -1.31%  "
Set-TextValue $ws.Range("D38") "2.878"
Set-TextValue $ws.Range("E38") "  +2.18%  "
Set-TextValue $ws.Range("D39") "1.225.72"
Set-TextValue $ws.Range("E39") "  +2.43%  "
Set-TextValue $ws.Range("D40") "0.01770"
Set-TextValue $ws.Range("E40") "  +0.94%  "
Set-TextValue $ws.Range("D41") "6.304"
Set-TextValue $ws.Range("E41") "  -1.22%  "
Set-TextValue $ws.Range("D42") "0.9174"
Set-TextValue $ws.Range("E42") "  +3.80%  "
Set-TextValue $ws.Range("E43") "  -0.90%  "
Set-TextValue $ws.Range("D44") "1.999.92"
Set-TextValue $ws.Range("E44") "  +1.64%  "
Set-TextValue $ws.Range("D45") "102.04"
Set-TextValue $ws.Range("E45") "  +1.08%  "
Set-TextValue $ws.Range("D46") "66.08"
Set-TextValue $ws.Range("E46") "  +3.00%  "
Set-TextValue $ws.Range("B47") "Mantle"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D47") "0.5082"
Set-TextValue $ws.Range("E47") "  -0.92%  "
Set-TextValue $ws.Range("B48") "BabyDogeCoin"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D48") "0.00000000120"
Set-TextValue $ws.Range("E48") "  -2.16%  "
Set-TextValue $ws.Range("D49") "9.281"
Set-TextValue $ws.Range("E49") "  +3.31%  "
Set-TextValue $ws.Range("D50") "0.4081"
Set-TextValue $ws.Range("E50") "  +2.58%  "
Set-TextValue $ws.Range("D51") "0.1140"
Set-TextValue $ws.Range("E51") "  +4.53%  "
